# Updates the cryptocurrency price/volume snapshot (and a few row re-orderings)
# to match the latest data refresh, per the commit:
# "Updated cryptos list on Thu Dec 14 02:42:33 UTC 2023 with GitHub Actions"
#
# Notes:
#  - Price (column D) values are plain text (they use '.' as a thousands
#    separator, e.g. "42.778.04"), so for values that Excel would otherwise
#    auto-coerce into a number (and possibly mangle, e.g. drop a trailing
#    zero or introduce floating point noise) we force the cell's
#    NumberFormat to Text ("@") before assigning the value.
#  - Volume(1h) (column E) values already contain padding spaces, so Excel
#    keeps them as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.778.04"
$ws.Range("E2").Value = "  +4.31%  "

$ws.Range("D3").Value = "2.248.98"
$ws.Range("E3").Value = "  +3.59%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.40"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.14"
$ws.Range("E7").Value = "  +4.72%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.659"
$ws.Range("E9").Value = "  +16.60%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.28"
$ws.Range("E10").Value = "  +9.75%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0971"
$ws.Range("E11").Value = "  +4.87%  "

$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "59.19"
$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.56"
$ws.Range("E13").Value = "  +9.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.104"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("D15").Value = "2.584.06"
$ws.Range("E15").Value = "  +3.81%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.885"
$ws.Range("E16").Value = "  +2.84%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.77"
$ws.Range("E17").Value = "  +4.29%  "

$ws.Range("D18").Value = "2.246.96"
$ws.Range("E18").Value = "  +3.02%  "

$ws.Range("D19").Value = "42.695.75"
$ws.Range("E19").Value = "  +4.56%  "

$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  +5.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("E21").Value = "  +3.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.02"
$ws.Range("E22").Value = "  +2.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.65"
$ws.Range("E23").Value = "  +2.58%  "

$ws.Range("E24").Value = "  -0.67%  "

$ws.Range("E25").Value = "  +5.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.55"
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.43"
$ws.Range("E28").Value = "  +0.38%  "

$ws.Range("E29").Value = "  -1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +10.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.50"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.96"
$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.61"
$ws.Range("E33").Value = "  +17.13%  "

$ws.Range("E34").Value = "  +5.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.40"
$ws.Range("E35").Value = "  +23.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0791"
$ws.Range("E36").Value = "  +6.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.125"
$ws.Range("E37").Value = "  +3.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  +7.56%  "

$ws.Range("E39").Value = "  +4.05%  "

$ws.Range("E40").Value = "  +7.25%  "

$ws.Range("E41").Value = "  +6.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.45"
$ws.Range("E42").Value = "  +7.51%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.81"
$ws.Range("E43").Value = "  +5.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.84"
$ws.Range("E44").Value = "  +2.34%  "

$ws.Range("E45").Value = "  +4.58%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.95"
$ws.Range("E46").Value = "  +5.13%  "

$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.81"
$ws.Range("E47").Value = "  +0.71%  "

$ws.Range("E48").Value = "  +5.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("E50").Value = "  -0.35%  "

$ws.Range("E51").Value = "  +4.28%  "
